# DailyTav.xlsx edit: add a new "AcctCode" field row to the DBD layout sheet
# and add a matching key-search definition row to the DBS sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "DBD" sheet
$ws2 = $wb.Worksheets.Item(2)   # "DBS" sheet

# --- Sheet "DBD": insert a new row 15 (AcctCode / 業務科目) --------------
# Old row 15 (CreateDate) and everything below it shifts down by one.
$ws1.Rows.Item(15).Insert(-4121, -4163)   # xlShiftDown, xlFormatFromLeftOrAbove

# Copy the formatting of the (now shifted) row 16 into the new blank row 15
# so the new row matches the rest of the table's look (thin borders etc.)
$ws1.Range("A16:G16").Copy()
$ws1.Range("A15:G15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row's content
$ws1.Range("A15").Value = 7
$ws1.Range("B15").Value = "AcctCode"
$ws1.Range("C15").Value = "業務科目"
$ws1.Range("D15").Value = "VARCHAR2"
$ws1.Range("E15").Value = 3
$ws1.Range("G15").Value = "TAV,TLD"

# The row-insert shifts the whole row (including the SEQ numbers in column A)
# down with it, but the SEQ column is really just the row's ordinal position,
# so re-stamp the correct sequential numbers for the rows that moved.
$ws1.Range("A16").Value = 8
$ws1.Range("A17").Value = 9
$ws1.Range("A18").Value = 10
$ws1.Range("A19").Value = 11

# Give the note cell (F15) its own distinct style: 細明體 12pt, centered,
# thin box border (matches the new cellXfs/font entries added to styles.xml)
$f15 = $ws1.Range("F15")
$f15.WrapText = $false
$f15.Font.Name = "細明體"
$f15.Font.Size = 12
$f15.HorizontalAlignment = -4108   # xlCenter
$f15.VerticalAlignment = -4108     # xlCenter

# --- Sheet "DBS": add a new key-search definition row --------------------
$ws2.Range("A2").Value = "CustNoAcDateRange"
$ws2.Range("B2").Value = "CustNo = ,AND AcDate >= ,AND AcDate <="
$ws2.Range("C2").Value = "AcDate asc"

# --- Restore selections / active sheet -----------------------------------
$ws2.Activate()
$ws2.Range("B5").Select()

$ws1.Activate()
$ws1.Range("G20").Select()
